$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HAL")

# Insert two new columns before D, shifting existing quarterly data (D:K) to (F:M)
$ws.Columns("D:E").Insert()

# Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30)
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 5936000
$ws.Cells.Item(8, 5).Value = 6172000
$ws.Cells.Item(9, 4).Value = 5269000
$ws.Cells.Item(9, 5).Value = 5384000
$ws.Cells.Item(10, 4).Value = 667000
$ws.Cells.Item(10, 5).Value = 788000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 5328000
$ws.Cells.Item(17, 5).Value = 5456000
$ws.Cells.Item(18, 4).Value = 608000
$ws.Cells.Item(18, 5).Value = 716000
$ws.Cells.Item(20, 4).Value = -2000
$ws.Cells.Item(20, 5).Value = -32000
$ws.Cells.Item(21, 4).Value = 1028000
$ws.Cells.Item(21, 5).Value = 1084000
$ws.Cells.Item(22, 4).Value = 148000
$ws.Cells.Item(22, 5).Value = 150000
$ws.Cells.Item(23, 4).Value = 458000
$ws.Cells.Item(23, 5).Value = 534000
$ws.Cells.Item(24, 4).Value = -163000
$ws.Cells.Item(24, 5).Value = 100000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 621000
$ws.Cells.Item(26, 5).Value = 434000
$ws.Cells.Item(27, 4).Value = 617000
$ws.Cells.Item(27, 5).Value = 435000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 47000
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 2000
$ws.Cells.Item(32, 5).Value = 32000
$ws.Cells.Item(33, 4).Value = 664000
$ws.Cells.Item(33, 5).Value = 435000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 664000
$ws.Cells.Item(35, 5).Value = 435000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 2008000
$ws.Cells.Item(41, 5).Value = 2057000
$ws.Cells.Item(42, 4).Value = "NA"
$ws.Cells.Item(42, 5).Value = "NA"
$ws.Cells.Item(43, 4).Value = 5234000
$ws.Cells.Item(43, 5).Value = 5526000
$ws.Cells.Item(44, 4).Value = 3028000
$ws.Cells.Item(44, 5).Value = 2887000
$ws.Cells.Item(45, 4).Value = 881000
$ws.Cells.Item(45, 5).Value = 966000
$ws.Cells.Item(46, 4).Value = 11151000
$ws.Cells.Item(46, 5).Value = 11436000
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 8961000
$ws.Cells.Item(48, 5).Value = 8821000
$ws.Cells.Item(49, 4).Value = 2825000
$ws.Cells.Item(49, 5).Value = 2800000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 3045000
$ws.Cells.Item(52, 5).Value = 2694000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 25982000
$ws.Cells.Item(54, 5).Value = 25751000
$ws.Cells.Item(57, 4).Value = 3018000
$ws.Cells.Item(57, 5).Value = 3142000
$ws.Cells.Item(58, 4).Value = 36000
$ws.Cells.Item(58, 5).Value = 35000
$ws.Cells.Item(59, 4).Value = 1748000
$ws.Cells.Item(59, 5).Value = 1776000
$ws.Cells.Item(60, 4).Value = 4802000
$ws.Cells.Item(60, 5).Value = 4953000
$ws.Cells.Item(61, 4).Value = 10421000
$ws.Cells.Item(61, 5).Value = 10424000
$ws.Cells.Item(62, 4).Value = 1215000
$ws.Cells.Item(62, 5).Value = 1357000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 16460000
$ws.Cells.Item(66, 5).Value = 16753000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 13739000
$ws.Cells.Item(72, 5).Value = 13216000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 9522000
$ws.Cells.Item(76, 5).Value = 8998000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 664000
$ws.Cells.Item(81, 5).Value = 435000
$ws.Cells.Item(83, 4).Value = 422000
$ws.Cells.Item(83, 5).Value = 400000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 852000
$ws.Cells.Item(89, 5).Value = 777000
$ws.Cells.Item(91, 4).Value = -551000
$ws.Cells.Item(91, 5).Value = -409000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -554000
$ws.Cells.Item(94, 5).Value = -2000
$ws.Cells.Item(96, 4).Value = -157000
$ws.Cells.Item(96, 5).Value = -157000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -338000
$ws.Cells.Item(100, 5).Value = -751000
$ws.Cells.Item(101, 4).Value = -9000
$ws.Cells.Item(101, 5).Value = -25000
$ws.Cells.Item(102, 4).Value = -49000
$ws.Cells.Item(102, 5).Value = -1000

# Minor restatements to previously reported Q3 2017 (now column I) figures
$ws.Cells.Item(9, 9).Value = 4747000
$ws.Cells.Item(10, 9).Value = 697000
$ws.Cells.Item(17, 9).Value = 4802000
$ws.Cells.Item(18, 9).Value = 642000
$ws.Cells.Item(20, 9).Value = -1000
$ws.Cells.Item(32, 9).Value = 1000
